$d = $word.ActiveDocument

$pairs = @(
  @("44×52=2288", "70×55=3850"),
  @("61×27=1647", "78×50=3900"),
  @("12×97=1164", "59×36=2124"),
  @("47×58=2726", "22×18=396"),
  @("89×62=5518", "79×32=2528"),
  @("38×89=3382", "73×87=6351"),
  @("43×74=3182", "12×31=372"),
  @("93×65=6045", "25×14=350"),
  @("40×61=2440", "48×28=1344"),
  @("36×14=504", "46×77=3542"),
  @("16×14=224", "64×36=2304"),
  @("20×78=1560", "63×49=3087"),
  @("78×34=2652", "90×54=4860"),
  @("37×17=629", "88×55=4840"),
  @("56×73=4088", "60×84=5040"),
  @("41×48=1968", "14×90=1260"),
  @("57×84=4788", "93×60=5580"),
  @("35×86=3010", "27×78=2106"),
  @("53×80=4240", "23×68=1564"),
  @("14×35=490", "19×66=1254"),
  @("79×56=4424", "49×40=1960"),
  @("17×62=1054", "33×73=2409"),
  @("31×20=620", "86×46=3956"),
  @("64×45=2880", "97×41=3977"),
  @("14×94=1316", "93×90=8370")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
